$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename header cells: "..._old" -> "..._FV2210", "..._new" -> "..._FV2304"
# ------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2210")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2304")
    }
}

# ------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1").
#    Creating a ListObject directly on top of the already-formatted
#    header row causes Excel to capture the existing header formatting
#    as a table "headerRowDxf" (and bloats styles.xml). To avoid that,
#    stash the header formatting on a scratch cell, strip the header's
#    direct formatting, build the table (with no table style), and
#    then restore the original header formatting from the scratch cell.
# ------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("AA1")

$headerRange.Cells.Item(1, 1).Copy()
$scratch.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$scratch.Clear()

# ------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
